$wb = $excel.ActiveWorkbook

# --- Add a new "State" column to the hotel_info sheet, right after "Hotel_Name" ---
$hotel = $wb.Worksheets.Item("hotel_info")
$hotel.Columns.Item(3).Insert()
$hotel.Cells.Item(1,3).Value = "State"
$hotel.Cells.Item(2,3).Value = "Louisiana"

# --- Reorder sheet tabs: review_info first, hotel_info second ---
$review = $wb.Worksheets.Item("review_info")
$review.Move($hotel)
